$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# ---------------------------------------------------------------------------
# 1) Build the row skeleton for the 10 new rows (231-240) by copying existing
#    rows that already carry the correct cell styles (s="2" dates, s="1" text
#    dates, s="4" percentage column). This keeps the style table untouched
#    (no new cellXfs are introduced) exactly like the source workbook.
#    Row 225 is an "open/pending" row (only C, D, L are populated) and is
#    used as the template for the two still-pending rows (231, 232).
#    Row 230 is a fully settled row and is used as the template for the
#    eight settled rows (233-240).
# ---------------------------------------------------------------------------
$ws.Range("A225:L225").Copy() | Out-Null
$ws.Range("A231:L231").Insert() | Out-Null

$ws.Range("A225:L225").Copy() | Out-Null
$ws.Range("A232:L232").Insert() | Out-Null

for ($i = 233; $i -le 240; $i++) {
    $ws.Range("A230:L230").Copy() | Out-Null
    $ws.Range("A" + $i + ":L" + $i).Insert() | Out-Null
}

$ws.Range("A231:B240").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 2) Dates (column C, serial numbers) and text labels (column D)
# ---------------------------------------------------------------------------
$dates = @{
    231 = 45287; 232 = 45288; 233 = 45289; 234 = 45289; 235 = 45289
    236 = 45289; 237 = 45289; 238 = 45290; 239 = 45290; 240 = 45291
}
$dLabels = @{
    231 = "2023-12-27"; 232 = "2023-12-28"; 233 = "2023-12-29"; 234 = "2023-12-29"
    235 = "2023-12-29"; 236 = "2023-12-29"; 237 = "2023-12-29"; 238 = "2023-12-30"
    239 = "2023-12-30"; 240 = "2023-12-31"
}
foreach ($r in 231..240) {
    $ws.Range("C" + $r).Value2 = $dates[$r]
    $ws.Range("D" + $r).Value2 = $dLabels[$r]
}

# ---------------------------------------------------------------------------
# 3) Rows 231 & 232 remain "pending" bets: only C/D/L are populated, so we
#    clear out everything the template row (225) might not already have
#    cleared (E, F, G, H, I, J, K).
# ---------------------------------------------------------------------------
$ws.Range("E231:K232").ClearContents() | Out-Null
$ws.Range("L231:L232").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 4) Rows 233-240: win/loss flag (E), odds (F), stake result (H), category
#    (J) and sub-category (K).
# ---------------------------------------------------------------------------
$ws.Range("E233").Value2 = 0
$ws.Range("E234").Value2 = 0
$ws.Range("E235").Value2 = 1
$ws.Range("E236").Value2 = 1
$ws.Range("E237").Value2 = 1
$ws.Range("E238").Value2 = 0
$ws.Range("E239").Value2 = 1
$ws.Range("E240").Value2 = 1

$ws.Range("F233").Value2 = 1.66
$ws.Range("F234").Value2 = 1.73
$ws.Range("F235").Value2 = 1.24
$ws.Range("F236").Value2 = 1.1
$ws.Range("F237").Value2 = 1.3
$ws.Range("F238").Value2 = 1.12
$ws.Range("F239").Value2 = 1.14
$ws.Range("F240").Value2 = 1.27

$ws.Range("H233").Value2 = -1405
$ws.Range("H234").Value2 = -8000
$ws.Range("H235").Value2 = 2400
$ws.Range("H236").Value2 = 500
$ws.Range("H237").Value2 = 2220
$ws.Range("H238").Value2 = -1120
$ws.Range("H239").Value2 = 560
$ws.Range("H240").Value2 = 1217

$ws.Range("J233").Value2 = "TENIS DE MESA"
$ws.Range("K233").Value2 = "SETKA CUP"

$ws.Range("J234").Value2 = "ESPORTS"
$ws.Range("K234").Value2 = "LOL EUROPEAN CIRCUIT"

$ws.Range("J235").Value2 = "ESPORTS"
$ws.Range("K235").Value2 = "LOL EUROPEAN CIRCUIT"

$ws.Range("J236").Value2 = "BASKET"
$ws.Range("K236").Value2 = "NBA"

$ws.Range("J237").Value2 = "BASKET"
$ws.Range("K237").Value2 = "NBA"

$ws.Range("J238").Value2 = "TENIS DE MESA"
$ws.Range("K238").Value2 = "MASTERS"

$ws.Range("J239").Value2 = "VOLLEY"
$ws.Range("K239").Value2 = "ITALIA SERIE A3"

$ws.Range("J240").Value2 = "ESPORTS"
$ws.Range("K240").Value2 = "DEMACIA CUP"

# ---------------------------------------------------------------------------
# 5) Running pool formulas (G = previous I, I = G + H) and the percentage
#    column (L). G235:G240 is entered as a single fill so the engine groups
#    it into one shared formula, matching how a user would drag-fill it.
# ---------------------------------------------------------------------------
$ws.Range("G233").Formula = "=I230"
$ws.Range("G234").Formula = "=I233"
$ws.Range("G235:G240").Formula = "=I234"

$ws.Range("I233:I240").Formula = "=G233+H233"

$ws.Range("L233:L240").Formula = "=ROUND((I233/`$G`$31-1)*100, 3)+`$L`$29"
